$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Insert two new blank rows at 3:4 on Sheet1 (everything below shifts down by 2)
$ws1.Rows("3:4").Insert()

# Sheet2: add two new rows of text values (stored as Text-formatted strings so the
# leading/trailing zeros survive), used by the new Sheet1 formulas below.
# NOTE: "0123.50" is registered in the shared-strings table before "0123" (matches
# target uniqueCount ordering), so write A3 first even though A2 is visually first.
$ws2.Cells.Item(3, 1).NumberFormat = "@"
$ws2.Cells.Item(3, 1).Value = "0123.50"

$ws2.Cells.Item(2, 1).NumberFormat = "@"
$ws2.Cells.Item(2, 1).Value = "0123"

# Sheet1 new row 3: label + formula pulling the rational-looking text from Sheet2!A2
$ws1.Cells.Item(3, 1).Value = "General"
$ws1.Cells.Item(3, 2).NumberFormat = "@"
$ws1.Cells.Item(3, 2).Formula = "=Sheet2!A2"
$ws1.Cells.Item(3, 3).Formula = "=Sheet2!`$A`$1"

# Sheet1 new row 4: label + formula pulling the rational-looking text from Sheet2!A3
$ws1.Cells.Item(4, 1).Value = "General"
$ws1.Cells.Item(4, 2).NumberFormat = "@"
$ws1.Cells.Item(4, 2).Formula = "=Sheet2!A3"
$ws1.Cells.Item(4, 3).Formula = "=Sheet2!`$A`$1"

# PageSetup on Sheet2 (new in the target) - touching it serializes a <pageSetup> node
$ps2 = $ws2.PageSetup
$ps2.PaperSize = 9
$ps2.Orientation = 1

# Selections: Sheet2!A3, then leave Sheet1 active with B2 selected
$ws2.Cells.Item(3, 1).Select()
$ws1.Activate()
$ws1.Cells.Item(2, 2).Select()
